$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part A: move the "_GoBack" bookmark from the very end of the document
#         (after the trailing comma of the T-SQL snippet) to right after
#         the title text " & SQL Job Installation".
# -----------------------------------------------------------------------

# Remove the bookmark currently sitting at the end of the document.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Find the end of the title text so we know where the bookmark should go.
$titleRng = $d.Content
$titleRng.Find.Execute("SQL Job Installation", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$titleEnd = $titleRng.End

# Adding a bookmark exactly at the end of a paragraph (just before the
# paragraph mark) is unreliable, so nudge it out of the way temporarily:
# type a placeholder character after the title, drop the bookmark right
# before that placeholder, then remove the placeholder again.
$placeholder = $d.Range($titleEnd, $titleEnd)
$placeholder.InsertAfter("X")

$bookmarkSpot = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null

$placeholderRange = $d.Range($titleEnd, $titleEnd + 1)
$placeholderRange.Text = ""

# -----------------------------------------------------------------------
# Part B: insert a new "SSIS Credential Name" bullet before the existing
#         "SSIS Proxy Name" bullet.
# -----------------------------------------------------------------------

$proxyRng = $d.Content
$proxyRng.Find.Execute("SSIS Proxy Name", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$proxyStart = $proxyRng.Start

# Create a new (empty) list paragraph above the "SSIS Proxy Name" bullet;
# it inherits the same ListParagraph / numbering formatting.
$splitPoint = $d.Range($proxyStart, $proxyStart)
$splitPoint.InsertParagraphBefore()

$boldLabel = "SSIS Credential Name"
$plainRest = ": Name of the Credential used to create the Proxy a SQL Agent Job will run under."

$newParaRng = $d.Range($proxyStart, $proxyStart)
$newParaRng.InsertBefore($boldLabel + $plainRest)

$labelRng = $d.Range($proxyStart, $proxyStart + $boldLabel.Length)
$labelRng.Bold = 1

# -----------------------------------------------------------------------
# Part C: split the description sentence on the "SSIS Proxy Name" bullet
#         and append a new sentence about Credential/Proxy creation.
# -----------------------------------------------------------------------

$descRng = $d.Content
$descRng.Find.Execute( `
    "The proxy account needs to be configured before the utility is run.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$descStart = $descRng.Start
$descEnd = $descRng.End

$afterRng = $d.Range($descEnd, $descEnd)
$afterRng.InsertAfter( `
    "If specified in the XML file, a Credential and Proxy Account can be created during deployment.")

$oldSentenceRng = $d.Range($descStart, $descEnd)
$oldSentenceRng.Delete()
